$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-06-06 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-06-07 Friday", 2) | Out-Null
$d.Content.Find.Execute("76×43=", $true, $false, $false, $false, $false, $true, 1, $false, "18×11=", 2) | Out-Null
$d.Content.Find.Execute("26×93=", $true, $false, $false, $false, $false, $true, 1, $false, "38×74=", 2) | Out-Null
$d.Content.Find.Execute("37×52=", $true, $false, $false, $false, $false, $true, 1, $false, "24×81=", 2) | Out-Null
$d.Content.Find.Execute("59×47=", $true, $false, $false, $false, $false, $true, 1, $false, "18×63=", 2) | Out-Null
$d.Content.Find.Execute("21×98=", $true, $false, $false, $false, $false, $true, 1, $false, "72×17=", 2) | Out-Null
$d.Content.Find.Execute("49×55=", $true, $false, $false, $false, $false, $true, 1, $false, "16×27=", 2) | Out-Null
$d.Content.Find.Execute("90×44=", $true, $false, $false, $false, $false, $true, 1, $false, "81×78=", 2) | Out-Null
$d.Content.Find.Execute("39×71=", $true, $false, $false, $false, $false, $true, 1, $false, "56×88=", 2) | Out-Null
$d.Content.Find.Execute("50×46=", $true, $false, $false, $false, $false, $true, 1, $false, "31×39=", 2) | Out-Null
$d.Content.Find.Execute("97×24=", $true, $false, $false, $false, $false, $true, 1, $false, "20×12=", 2) | Out-Null
$d.Content.Find.Execute("71×64=", $true, $false, $false, $false, $false, $true, 1, $false, "98×87=", 2) | Out-Null
$d.Content.Find.Execute("61×91=", $true, $false, $false, $false, $false, $true, 1, $false, "22×69=", 2) | Out-Null
$d.Content.Find.Execute("92×40=", $true, $false, $false, $false, $false, $true, 1, $false, "80×90=", 2) | Out-Null
$d.Content.Find.Execute("62×75=", $true, $false, $false, $false, $false, $true, 1, $false, "53×21=", 2) | Out-Null
$d.Content.Find.Execute("87×77=", $true, $false, $false, $false, $false, $true, 1, $false, "13×73=", 2) | Out-Null
$d.Content.Find.Execute("86×16=", $true, $false, $false, $false, $false, $true, 1, $false, "67×42=", 2) | Out-Null
$d.Content.Find.Execute("77×36=", $true, $false, $false, $false, $false, $true, 1, $false, "58×93=", 2) | Out-Null
$d.Content.Find.Execute("13×87=", $true, $false, $false, $false, $false, $true, 1, $false, "88×31=", 2) | Out-Null
$d.Content.Find.Execute("87×14=", $true, $false, $false, $false, $false, $true, 1, $false, "51×33=", 2) | Out-Null
$d.Content.Find.Execute("79×59=", $true, $false, $false, $false, $false, $true, 1, $false, "73×26=", 2) | Out-Null
$d.Content.Find.Execute("46×62=", $true, $false, $false, $false, $false, $true, 1, $false, "77×96=", 2) | Out-Null
$d.Content.Find.Execute("29×43=", $true, $false, $false, $false, $false, $true, 1, $false, "95×59=", 2) | Out-Null
$d.Content.Find.Execute("16×75=", $true, $false, $false, $false, $false, $true, 1, $false, "77×30=", 2) | Out-Null
$d.Content.Find.Execute("24×22=", $true, $false, $false, $false, $false, $true, 1, $false, "52×14=", 2) | Out-Null
$d.Content.Find.Execute("19×71=", $true, $false, $false, $false, $false, $true, 1, $false, "65×80=", 2) | Out-Null
